$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 3), mirroring the values already present in row 2.
$rng = $ws.Range("A3:H3")

# Force the whole row to Text format first so that numeric-looking values
# ("333", "1111") are stored as text instead of being auto-converted to
# numbers, matching the original inlineStr/text cells used elsewhere in
# this sheet.
$rng.NumberFormat = "@"

$ws.Range("A3").Value = "address"
$ws.Range("B3").Value = "Ho Chi Minh"
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = "first"
$ws.Range("E3").Value = "last"
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = "333"
$ws.Range("H3").Value = "1111"

# Restore default (unstyled) cell formatting so the new row matches the
# plain/no-style formatting of row 2.
$rng.Style = "Normal"
